$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank separator row immediately above each section header
# (original, pre-edit row numbers; applied bottom-to-top so earlier
# insertion points keep their meaning).
$insertPositions = @(50, 46, 37, 35, 31, 26, 23, 20, 17, 13, 10, 7)
foreach ($p in $insertPositions) {
    $ws.Rows.Item($p).Insert()
}

# After the inserts above, the "numberformat" example rows (originally
# B38:B45 holding "123456" with a number-format style) now live at
# B48:B55. Move the formatted value into column C and replace column B
# with a descriptive label in the default style.
$labels = @("numberformat=currency", "numberformat=date", "numberformat=date_time", "numberformat=number", "numberformat=percent", "numberformat=text", "numberformat=time", "numberformat=scientific")
$row = 48
foreach ($label in $labels) {
    $src = $ws.Cells.Item($row, 2)
    $dst = $ws.Cells.Item($row, 3)
    $src.Copy($dst)
    $src.ClearFormats()
    $src.Value = $label
    $row = $row + 1
}
